# Auto-generated Excel COM-interop script applying the Maduin_Profits diff
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets of the workbook.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40 (Leve Item ID 5505)
$ws.Cells.Item(40, 8).Value = 2977.6
$ws.Cells.Item(40, 9).Value = 2944
$ws.Cells.Item(40, 11).Value = 2944
$ws.Cells.Item(40, 13).Value = -2769

# Row 64 (Leve Item ID 5506)
$ws.Cells.Item(64, 8).Value = 11511.111
$ws.Cells.Item(64, 9).Value = 4900
$ws.Cells.Item(64, 10).Value = 16800
$ws.Cells.Item(64, 11).Value = 4900
$ws.Cells.Item(64, 12).Value = 16800
$ws.Cells.Item(64, 13).Value = -4652
$ws.Cells.Item(64, 14).Value = -17296

# Row 67 (Leve Item ID 5506)
$ws.Cells.Item(67, 8).Value = 11511.111
$ws.Cells.Item(67, 9).Value = 4900
$ws.Cells.Item(67, 10).Value = 16800
$ws.Cells.Item(67, 11).Value = 4900
$ws.Cells.Item(67, 12).Value = 16800
$ws.Cells.Item(67, 13).Value = -4042
$ws.Cells.Item(67, 14).Value = -18516

# Row 74 (Leve Item ID 5507)
$ws.Cells.Item(74, 8).Value = 3201
$ws.Cells.Item(74, 9).Value = 2801.5
$ws.Cells.Item(74, 10).Value = 4000
$ws.Cells.Item(74, 11).Value = 2801.5
$ws.Cells.Item(74, 12).Value = 4000
$ws.Cells.Item(74, 13).Value = -1865.5
$ws.Cells.Item(74, 14).Value = -5872

# Row 76 (Leve Item ID 12602)
$ws.Cells.Item(76, 9).Value = 0
$ws.Cells.Item(76, 11).Value = 0
$ws.Cells.Item(76, 13).ClearContents()

# Row 77 (Leve Item ID 5507)
$ws.Cells.Item(77, 8).Value = 3201
$ws.Cells.Item(77, 9).Value = 2801.5
$ws.Cells.Item(77, 10).Value = 4000
$ws.Cells.Item(77, 11).Value = 14007.5
$ws.Cells.Item(77, 12).Value = 20000
$ws.Cells.Item(77, 13).Value = -9327.5
$ws.Cells.Item(77, 14).Value = -29360

# Row 79 (Leve Item ID 12602)
$ws.Cells.Item(79, 9).Value = 0
$ws.Cells.Item(79, 11).Value = 0
$ws.Cells.Item(79, 13).ClearContents()

# Row 95 (Leve Item ID 18200)
$ws.Cells.Item(95, 8).Value = 17750
$ws.Cells.Item(95, 10).Value = 17750
$ws.Cells.Item(95, 12).Value = 17750
$ws.Cells.Item(95, 14).Value = -23242

# Row 106 (Leve Item ID 19903)
$ws.Cells.Item(106, 8).Value = 5000
$ws.Cells.Item(106, 9).Value = 5000
$ws.Cells.Item(106, 11).Value = 5000
$ws.Cells.Item(106, 13).Value = -4369

$ws = $wb.Worksheets.Item("ARM")
# Row 2 (Leve Item ID 27713)
$ws.Cells.Item(2, 8).Value = 2383.3333
$ws.Cells.Item(2, 9).Value = 857.6
$ws.Cells.Item(2, 11).Value = 857.6
$ws.Cells.Item(2, 13).Value = -744.6

# Row 32 (Leve Item ID 44147)
$ws.Cells.Item(32, 8).Value = 2288.4595
$ws.Cells.Item(32, 9).Value = 2288.4595
$ws.Cells.Item(32, 11).Value = 2288.4595
$ws.Cells.Item(32, 13).Value = -2001.4595

# Row 63 (Leve Item ID 12528)
$ws.Cells.Item(63, 8).Value = 6312.5
$ws.Cells.Item(63, 10).Value = 7916.6665
$ws.Cells.Item(63, 12).Value = 7916.6665
$ws.Cells.Item(63, 14).Value = -9288.666499999999

# Row 66 (Leve Item ID 12528)
$ws.Cells.Item(66, 8).Value = 6312.5
$ws.Cells.Item(66, 10).Value = 7916.6665
$ws.Cells.Item(66, 12).Value = 39583.3325
$ws.Cells.Item(66, 14).Value = -46447.3325

# Row 116 (Leve Item ID 27713)
$ws.Cells.Item(116, 8).Value = 2383.3333
$ws.Cells.Item(116, 9).Value = 857.6
$ws.Cells.Item(116, 11).Value = 857.6
$ws.Cells.Item(116, 13).Value = 1436.4

# Row 122 (Leve Item ID 36168)
$ws.Cells.Item(122, 8).Value = 1553.0834
$ws.Cells.Item(122, 9).Value = 1526.5
$ws.Cells.Item(122, 10).Value = 1606.25
$ws.Cells.Item(122, 11).Value = 4579.5
$ws.Cells.Item(122, 12).Value = 4818.75
$ws.Cells.Item(122, 13).Value = -2129.5
$ws.Cells.Item(122, 14).Value = -9718.75

# Row 135 (Leve Item ID 42016)
$ws.Cells.Item(135, 8).Value = 8379499
$ws.Cells.Item(135, 9).Value = 50000000
$ws.Cells.Item(135, 10).Value = 55398.8
$ws.Cells.Item(135, 11).Value = 50000000
$ws.Cells.Item(135, 12).Value = 55398.8
$ws.Cells.Item(135, 13).Value = -49994930
$ws.Cells.Item(135, 14).Value = -65538.8

$ws = $wb.Worksheets.Item("BSM")
# Row 3 (Leve Item ID 27713)
$ws.Cells.Item(3, 8).Value = 2383.3333
$ws.Cells.Item(3, 9).Value = 857.6
$ws.Cells.Item(3, 11).Value = 857.6
$ws.Cells.Item(3, 13).Value = -743.6

# Row 82 (Leve Item ID 11877)
$ws.Cells.Item(82, 8).Value = 16000
$ws.Cells.Item(82, 10).Value = 0
$ws.Cells.Item(82, 12).Value = 0
$ws.Cells.Item(82, 14).ClearContents()

# Row 85 (Leve Item ID 11877)
$ws.Cells.Item(85, 8).Value = 16000
$ws.Cells.Item(85, 10).Value = 0
$ws.Cells.Item(85, 12).Value = 0
$ws.Cells.Item(85, 14).ClearContents()

# Row 86 (Leve Item ID 12526)
$ws.Cells.Item(86, 8).Value = 2022.091
$ws.Cells.Item(86, 9).Value = 2068.8
$ws.Cells.Item(86, 10).Value = 1983.1666
$ws.Cells.Item(86, 11).Value = 2068.8
$ws.Cells.Item(86, 12).Value = 1983.1666
$ws.Cells.Item(86, 13).Value = -945.8000000000002
$ws.Cells.Item(86, 14).Value = -4229.1666

# Row 89 (Leve Item ID 12526)
$ws.Cells.Item(89, 8).Value = 2022.091
$ws.Cells.Item(89, 9).Value = 2068.8
$ws.Cells.Item(89, 10).Value = 1983.1666
$ws.Cells.Item(89, 11).Value = 10344
$ws.Cells.Item(89, 12).Value = 9915.833000000001
$ws.Cells.Item(89, 13).Value = -4728
$ws.Cells.Item(89, 14).Value = -21147.833

# Row 95 (Leve Item ID 18194)
$ws.Cells.Item(95, 8).Value = 3016
$ws.Cells.Item(95, 10).Value = 3016
$ws.Cells.Item(95, 12).Value = 3016
$ws.Cells.Item(95, 14).Value = -8508

# Row 103 (Leve Item ID 18514)
$ws.Cells.Item(103, 8).Value = 25000
$ws.Cells.Item(103, 10).Value = 25000
$ws.Cells.Item(103, 12).Value = 25000
$ws.Cells.Item(103, 14).Value = -27344

# Row 105 (Leve Item ID 19947)
$ws.Cells.Item(105, 8).Value = 3121.3333
$ws.Cells.Item(105, 9).Value = 2545.6
$ws.Cells.Item(105, 11).Value = 2545.6
$ws.Cells.Item(105, 13).Value = -798.5999999999999

$ws = $wb.Worksheets.Item("CRP")
# Row 22 (Leve Item ID 5367)
$ws.Cells.Item(22, 8).Value = 1443.3334
$ws.Cells.Item(22, 9).Value = 998
$ws.Cells.Item(22, 11).Value = 998
$ws.Cells.Item(22, 13).Value = -648

# Row 92 (Leve Item ID 18041)
$ws.Cells.Item(92, 8).Value = 39998.5
$ws.Cells.Item(92, 10).Value = 39998.5
$ws.Cells.Item(92, 12).Value = 39998.5
$ws.Cells.Item(92, 14).Value = -44990.5

# Row 96 (Leve Item ID 18193)
$ws.Cells.Item(96, 8).Value = 14249.25
$ws.Cells.Item(96, 10).Value = 14249.25
$ws.Cells.Item(96, 12).Value = 14249.25
$ws.Cells.Item(96, 14).Value = -19741.25

# Row 106 (Leve Item ID 18661)
$ws.Cells.Item(106, 8).Value = 189833.33
$ws.Cells.Item(106, 10).Value = 189833.33
$ws.Cells.Item(106, 12).Value = 189833.33
$ws.Cells.Item(106, 14).Value = -192357.33

# Row 132 (Leve Item ID 44019)
$ws.Cells.Item(132, 8).Value = 3703
$ws.Cells.Item(132, 9).Value = 3982.6667
$ws.Cells.Item(132, 11).Value = 11948.0001
$ws.Cells.Item(132, 13).Value = -9418.000100000001

$ws = $wb.Worksheets.Item("CUL")
# Row 131 (Leve Item ID 36060)
$ws.Cells.Item(131, 8).Value = 892.03845
$ws.Cells.Item(131, 10).Value = 943.3913
$ws.Cells.Item(131, 12).Value = 2830.1739
$ws.Cells.Item(131, 14).Value = -12910.1739

# Row 140 (Leve Item ID 44097)
$ws.Cells.Item(140, 8).Value = 2302.875
$ws.Cells.Item(140, 9).Value = 693
$ws.Cells.Item(140, 10).Value = 4986
$ws.Cells.Item(140, 11).Value = 2079
$ws.Cells.Item(140, 12).Value = 14958
$ws.Cells.Item(140, 13).Value = 3101
$ws.Cells.Item(140, 14).Value = -25318

$ws = $wb.Worksheets.Item("GSM")
# Row 2 (Leve Item ID 5062)
$ws.Cells.Item(2, 8).Value = 881.8333
$ws.Cells.Item(2, 9).Value = 1223.9166
$ws.Cells.Item(2, 10).Value = 197.66667
$ws.Cells.Item(2, 11).Value = 1223.9166
$ws.Cells.Item(2, 12).Value = 197.66667
$ws.Cells.Item(2, 13).Value = -1110.9166
$ws.Cells.Item(2, 14).Value = -423.66667

# Row 43 (Leve Item ID 4218)
$ws.Cells.Item(43, 8).Value = 27999
$ws.Cells.Item(43, 10).Value = 54999
$ws.Cells.Item(43, 12).Value = 54999
$ws.Cells.Item(43, 14).Value = -55301

# Row 47 (Leve Item ID 4343)
$ws.Cells.Item(47, 8).Value = 0
$ws.Cells.Item(47, 10).Value = 0
$ws.Cells.Item(47, 12).Value = 0
$ws.Cells.Item(47, 14).ClearContents()

# Row 126 (Leve Item ID 36184)
$ws.Cells.Item(126, 8).Value = 8724.75
$ws.Cells.Item(126, 9).Value = 8724.75
$ws.Cells.Item(126, 11).Value = 26174.25
$ws.Cells.Item(126, 13).Value = -23704.25

# Row 136 (Leve Item ID 42218)
$ws.Cells.Item(136, 8).Value = 28345.223
$ws.Cells.Item(136, 10).Value = 28345.223
$ws.Cells.Item(136, 12).Value = 85035.66900000001
$ws.Cells.Item(136, 14).Value = -90135.66900000001

$ws = $wb.Worksheets.Item("LTW")
# Row 16 (Leve Item ID 5289)
$ws.Cells.Item(16, 8).Value = 3530.3333
$ws.Cells.Item(16, 9).Value = 295.5
$ws.Cells.Item(16, 10).Value = 10000
$ws.Cells.Item(16, 11).Value = 295.5
$ws.Cells.Item(16, 12).Value = 10000
$ws.Cells.Item(16, 13).Value = -125.5
$ws.Cells.Item(16, 14).Value = -10340

# Row 22 (Leve Item ID 5277)
$ws.Cells.Item(22, 8).Value = 1248.625
$ws.Cells.Item(22, 9).Value = 1509.75
$ws.Cells.Item(22, 10).Value = 987.5
$ws.Cells.Item(22, 11).Value = 1509.75
$ws.Cells.Item(22, 12).Value = 987.5
$ws.Cells.Item(22, 13).Value = -1214.75
$ws.Cells.Item(22, 14).Value = -1577.5

# Row 27 (Leve Item ID 5277)
$ws.Cells.Item(27, 8).Value = 1248.625
$ws.Cells.Item(27, 9).Value = 1509.75
$ws.Cells.Item(27, 10).Value = 987.5
$ws.Cells.Item(27, 11).Value = 1509.75
$ws.Cells.Item(27, 12).Value = 987.5
$ws.Cells.Item(27, 13).Value = -1402.75
$ws.Cells.Item(27, 14).Value = -1201.5
